# Rename the C#-style `WriteLine(` call to the Python-style `write_line(`
# inside the "    WriteLine(" code-sample run that appears once per slide
# (inside a grouped code-block textbox) across all 16 slides.

function Find-TextFrameShape($shapes, $needle) {
    # Recursively search a Shapes/GroupShapes collection (including nested
    # groups) for the first shape whose text contains $needle.
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Type -eq 6) {
            # msoGroup -> recurse into the group's children
            $found = Find-TextFrameShape $sh.GroupItems $needle
            if ($found -ne $null) {
                return $found
            }
        } elseif ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -match [regex]::Escape($needle)) {
                return $sh
            }
        }
    }
    return $null
}

$oldText = "    WriteLine("
$newText = "    write_line("

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    $shape = Find-TextFrameShape $slide.Shapes $oldText
    if ($shape -eq $null) {
        continue
    }

    $tr = $shape.TextFrame.TextRange
    $full = $tr.Text
    $idx = $full.IndexOf($oldText)
    if ($idx -lt 0) {
        continue
    }

    # Grab just the single run's worth of characters ("    WriteLine(") and
    # overwrite its text in place, leaving every other run (formatting,
    # the following "i" / ");" runs, etc.) untouched.
    $sub = $tr.Characters($idx + 1, $oldText.Length)
    $sub.Text = $newText
}
